# Apply the edit described by the commit: re-sort the correlation table by
# the "origXxx field" column (B) instead of the previous "avgCreditScore"
# correlation column (E), widen the AutoFilter / _FilterDatabase range to
# start at column A instead of column B, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Sort the data rows (2:17) by column B, ascending, header excluded ---
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("B1:B17"))
$sortObj.SetRange($ws.Range("A1:R17"))
$sortObj.Header = 1
$sortObj.Apply()

# --- 2. Re-apply the AutoFilter so it spans A1:R17 (was B1:R17) ---
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:R17").AutoFilter()

# --- 3. Update the hidden _FilterDatabase defined name to match ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "corrFMAcqFMac!_FilterDatabase") {
        $n.RefersTo = "=corrFMAcqFMac!`$A`$1:`$R`$17"
    }
}

# --- 4. Move the active selection to A13:A15 ---
$ws.Range("A13:A15").Select()
